$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.5
$ws.Range("H2").Value = 4.1
$ws.Range("I2").Value = 6.5
$ws.Range("J2").Value = 2.1
$ws.Range("L2").Value = 7
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.65
$ws.Range("U2").Value = 2.38
$ws.Range("V2").Value = 1.53
$ws.Range("W2").Value = 5
$ws.Range("Z2").Value = 10
$ws.Range("AB2").Value = 41
$ws.Range("AF2").Value = 101
$ws.Range("AH2").Value = 13
$ws.Range("AI2").Value = 34
$ws.Range("AJ2").Value = 21
$ws.Range("AK2").Value = 81
$ws.Range("AM2").Value = 67
$ws.Range("AW2").Value = 8
$ws.Range("AX2").Value = 41
$ws.Range("J4").Value = 3.35
$ws.Range("K4").Value = 2.07
$ws.Range("L4").Value = 2.95
$ws.Range("O4").Value = 1.3
$ws.Range("V4").Value = 1.98
$ws.Range("W4").Value = 8.75
$ws.Range("X4").Value = 14.5
$ws.Range("AA4").Value = 24
$ws.Range("AB4").Value = 32
$ws.Range("AH4").Value = 8.75
$ws.Range("AI4").Value = 13
$ws.Range("AM4").Value = 27
$ws.Range("AO4").Value = 15
$ws.Range("AP4").Value = 20
$ws.Range("AR4").Value = 90
$ws.Range("AS4").Value = 250
$ws.Range("AT4").Value = 2.62
$ws.Range("AV4").Value = 50
$ws.Range("AX4").Value = 12.5
$ws.Range("AY4").Value = 17.5
$ws.Range("AZ4").Value = 50
$ws.Range("BA4").Value = 70
$ws.Range("BB4").Value = 175
